# Update Sheet1 statistics cells (ログイン者(人) / 入力者(人) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 27
$ws.Range("C3").Value = 47
$ws.Range("D3").Value = 40
$ws.Range("D4").Value = 69
$ws.Range("D6").Value = 89
$ws.Range("C7").Value = 77
$ws.Range("D8").Value = 33
$ws.Range("C9").Value = 53
$ws.Range("D9").Value = 45
$ws.Range("C10").Value = 30
$ws.Range("C11").Value = 66
$ws.Range("D11").Value = 56
$ws.Range("C12").Value = 40
$ws.Range("D12").Value = 33
$ws.Range("C13").Value = 51
$ws.Range("D13").Value = 37
$ws.Range("C14").Value = 98
$ws.Range("D14").Value = 87
$ws.Range("C16").Value = 76
$ws.Range("D16").Value = 61
$ws.Range("C17").Value = 38
$ws.Range("D17").Value = 31
$ws.Range("C21").Value = 67
$ws.Range("D21").Value = 53
$ws.Range("D22").Value = 28
$ws.Range("C23").Value = 27
$ws.Range("D23").Value = 22
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 33
$ws.Range("C26").Value = 41
$ws.Range("D26").Value = 34
$ws.Range("C27").Value = 61
$ws.Range("D27").Value = 46
$ws.Range("C28").Value = 59
$ws.Range("D28").Value = 50
$ws.Range("C29").Value = 70
$ws.Range("D29").Value = 55
$ws.Range("C30").Value = 57
$ws.Range("D30").Value = 46
$ws.Range("C31").Value = 63
$ws.Range("D31").Value = 49
$ws.Range("C33").Value = 54
$ws.Range("D33").Value = 41
$ws.Range("D34").Value = 45
$ws.Range("C35").Value = 61
$ws.Range("D35").Value = 53
$ws.Range("D36").Value = 28
$ws.Range("C37").Value = 61
$ws.Range("D37").Value = 44
$ws.Range("C40").Value = 83
$ws.Range("D40").Value = 73
$ws.Range("C42").Value = 77
$ws.Range("D42").Value = 65
$ws.Range("C43").Value = 65
$ws.Range("D43").Value = 56
$ws.Range("C44").Value = 70
$ws.Range("D44").Value = 59
$ws.Range("C46").Value = 72
$ws.Range("D46").Value = 57
$ws.Range("C47").Value = 56
$ws.Range("D47").Value = 48
$ws.Range("C48").Value = 53
$ws.Range("C49").Value = 63
$ws.Range("D49").Value = 54
$ws.Range("C52").Value = 53
$ws.Range("C55").Value = 66
$ws.Range("D55").Value = 47
$ws.Range("C57").Value = 73
$ws.Range("D57").Value = 46
$ws.Range("C58").Value = 75
$ws.Range("D58").Value = 60
$ws.Range("C60").Value = 58
$ws.Range("D60").Value = 40
$ws.Range("C61").Value = 41
$ws.Range("D61").Value = 35
$ws.Range("C62").Value = 64
$ws.Range("D62").Value = 54
$ws.Range("C63").Value = 108
$ws.Range("D63").Value = 100
$ws.Range("D65").Value = 28
$ws.Range("C66").Value = 67
$ws.Range("D66").Value = 50
$ws.Range("C68").Value = 62
$ws.Range("D68").Value = 50
$ws.Range("D70").Value = 45
$ws.Range("C72").Value = 43
$ws.Range("C73").Value = 61
$ws.Range("D73").Value = 41
$ws.Range("D75").Value = 106
$ws.Range("C76").Value = 29
$ws.Range("D76").Value = 15
$ws.Range("C78").Value = 82
$ws.Range("D78").Value = 64
$ws.Range("C81").Value = 77
$ws.Range("D81").Value = 55
$ws.Range("C82").Value = 30
$ws.Range("D82").Value = 16
$ws.Range("D83").Value = 64
$ws.Range("C92").Value = 266
$ws.Range("D92").Value = 193
$ws.Range("C93").Value = 5616
$ws.Range("D93").Value = 4536

# Add two new empty worksheets after Sheet1, named Sheet2 and Sheet3
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# Re-select Sheet1 as the active sheet
$ws.Select()
